# Apply "added city meta data" change:
# Set "X" markers in the Meta data columns (N, P, Q, R) for several city rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 3, 4, 5: N, P, R get "X"
foreach ($r in 3, 4, 5) {
    $ws.Range("N$r").Value = "X"
    $ws.Range("P$r").Value = "X"
    $ws.Range("R$r").Value = "X"
}

# Rows 6, 7, 8, 9, 10, 12: N, R get "X"
foreach ($r in 6, 7, 8, 9, 10, 12) {
    $ws.Range("N$r").Value = "X"
    $ws.Range("R$r").Value = "X"
}

# Update the active selection to T10 (matches post-edit sheetView selection)
$ws.Range("T10").Select()

# Best-effort: reposition the window to match the recorded workbook view state.
$win = $excel.ActiveWindow
$win.Left = 47200
$win.Top = -8820
